$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 191, shifting existing rows 191..256 down to 192..257
$ws.Rows(191).Insert()

# Populate the newly inserted row 191 with its data
$ws.Range("A191").Value = 4
$ws.Range("B191").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C191").Value = "Los Lagos"
$ws.Range("D191").Value = 44809
$ws.Range("E191").Value = 10
$ws.Range("F191").Value = 100112039
$ws.Range("G191").Value = "Ciboulette"
$ws.Range("H191").Value = "Sin especificar"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 80
$ws.Range("K191").Value = 3500
$ws.Range("L191").Value = 3500
$ws.Range("M191").Value = 3500
$ws.Range("N191").Value = "`$/docena de atados"
$ws.Range("O191").Value = "Región Metropolitana"
$ws.Range("P191").Value = 1167
$ws.Range("Q191").Value = 3
$ws.Range("R191").Value = "Hortaliza"
